$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$fmtSource = $ws.Range("H2")

$ws.Range("E2").Value = "2026-02-28 17:18:13"
$ws.Range("K2").Value = "11.7 MJ/m2"
$ws.Range("E3").Value = "2026-02-28 17:18:15"
$ws.Range("K3").Value = "11.3 MJ/m2"
$ws.Range("E4").Value = "2026-02-28 17:18:18"
$ws.Range("J4").Value = "1024.7 hPa"
$ws.Range("K4").Value = "6.0 MJ/m2"
$ws.Range("O4").Value = "11.1 °C"
$ws.Range("E5").Value = "2026-02-28 17:18:20"
$ws.Range("I5").Value = "0.7 mm"
$ws.Range("N5").Value = "-2.3 °C 16:44 TU"
$ws.Range("E6").Value = "2026-02-28 17:18:22"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "83%"
$fmtSource.Copy()
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("J6").Value = "1024.6 hPa"
$ws.Range("K6").Value = "10.3 MJ/m2"
$ws.Range("L6").Value = "23.0 km/h - 36º 16:52 TU"
$ws.Range("E7").Value = "2026-02-28 17:18:25"
$ws.Range("K7").Value = "5.8 MJ/m2"
$ws.Range("E8").Value = "2026-02-28 17:18:27"
$ws.Range("K8").Value = "4.3 MJ/m2"
$ws.Range("E9").Value = "2026-02-28 17:18:29"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "82%"
$fmtSource.Copy()
$ws.Range("H9").PasteSpecial(-4122)
$ws.Range("K9").Value = "11.3 MJ/m2"
$ws.Range("L9").Value = "28.4 km/h - 312º 16:39 TU"
$ws.Range("O9").Value = "11.3 °C"
$ws.Range("E10").Value = "2026-02-28 17:18:32"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "85%"
$fmtSource.Copy()
$ws.Range("H10").PasteSpecial(-4122)
$ws.Range("K10").Value = "8.7 MJ/m2"
$ws.Range("O10").Value = "11.1 °C"
$ws.Range("E11").Value = "2026-02-28 17:18:34"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "72%"
$fmtSource.Copy()
$ws.Range("H11").PasteSpecial(-4122)
$ws.Range("O11").Value = "6.9 °C"
$ws.Range("E12").Value = "2026-02-28 17:18:37"
$ws.Range("O12").Value = "10.7 °C"
$ws.Range("E13").Value = "2026-02-28 17:18:39"
$ws.Range("E14").Value = "2026-02-28 17:18:41"
$ws.Range("K14").Value = "5.7 MJ/m2"
$ws.Range("E15").Value = "2026-02-28 17:18:44"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "80%"
$fmtSource.Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("O15").Value = "11.2 °C"
$ws.Range("E16").Value = "2026-02-28 17:18:46"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "62%"
$fmtSource.Copy()
$ws.Range("H16").PasteSpecial(-4122)
$ws.Range("K16").Value = "11.5 MJ/m2"
$ws.Range("E17").Value = "2026-02-28 17:18:48"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "79%"
$fmtSource.Copy()
$ws.Range("H17").PasteSpecial(-4122)
$ws.Range("E18").Value = "2026-02-28 17:18:51"
$ws.Range("K18").Value = "9.7 MJ/m2"
$ws.Range("O18").Value = "11.8 °C"
$ws.Range("E19").Value = "2026-02-28 17:18:53"
$ws.Range("K19").Value = "7.7 MJ/m2"
$ws.Range("O19").Value = "8.0 °C"
$ws.Range("E20").Value = "2026-02-28 17:18:56"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "56%"
$fmtSource.Copy()
$ws.Range("H20").PasteSpecial(-4122)
$ws.Range("K20").Value = "13.9 MJ/m2"
$ws.Range("E21").Value = "2026-02-28 17:18:58"
$ws.Range("J21").Value = "1024.0 hPa"
$ws.Range("K21").Value = "10.9 MJ/m2"
$ws.Range("O21").Value = "7.4 °C"
$ws.Range("E22").Value = "2026-02-28 17:19:01"
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "63%"
$fmtSource.Copy()
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("K22").Value = "11.4 MJ/m2"
$ws.Range("E23").Value = "2026-02-28 17:19:03"
$ws.Range("K23").Value = "12.2 MJ/m2"
$ws.Range("O23").Value = "-0.1 °C"
$ws.Range("E24").Value = "2026-02-28 17:19:05"
$ws.Range("K24").Value = "2.7 MJ/m2"
$ws.Range("E25").Value = "2026-02-28 17:19:08"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = "56%"
$fmtSource.Copy()
$ws.Range("H25").PasteSpecial(-4122)
$ws.Range("O25").Value = "1.6 °C"
$ws.Range("E26").Value = "2026-02-28 17:19:10"
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = "77%"
$fmtSource.Copy()
$ws.Range("H26").PasteSpecial(-4122)
$ws.Range("E27").Value = "2026-02-28 17:19:12"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = "48%"
$fmtSource.Copy()
$ws.Range("H27").PasteSpecial(-4122)
$ws.Range("K27").Value = "13.8 MJ/m2"
$ws.Range("E28").Value = "2026-02-28 17:19:15"
$ws.Range("J28").Value = "1024.7 hPa"
$ws.Range("K28").Value = "7.1 MJ/m2"
$ws.Range("O28").Value = "9.4 °C"
$ws.Range("E29").Value = "2026-02-28 17:19:17"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "84%"
$fmtSource.Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("K29").Value = "12.3 MJ/m2"
$ws.Range("E30").Value = "2026-02-28 17:19:20"
$ws.Range("K30").Value = "12.6 MJ/m2"
$ws.Range("E31").Value = "2026-02-28 17:19:22"
$ws.Range("K31").Value = "11.9 MJ/m2"
$ws.Range("L31").Value = "59.0 km/h - 356º 16:55 TU"
$ws.Range("E32").Value = "2026-02-28 17:19:24"
$ws.Range("E33").Value = "2026-02-28 17:19:27"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "63%"
$fmtSource.Copy()
$ws.Range("H33").PasteSpecial(-4122)
$ws.Range("K33").Value = "11.7 MJ/m2"
$ws.Range("E34").Value = "2026-02-28 17:19:29"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "62%"
$fmtSource.Copy()
$ws.Range("H34").PasteSpecial(-4122)
$ws.Range("E35").Value = "2026-02-28 17:19:32"
$ws.Range("E36").Value = "2026-02-28 17:19:34"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "83%"
$fmtSource.Copy()
$ws.Range("H36").PasteSpecial(-4122)
$ws.Range("K36").Value = "13.5 MJ/m2"
$ws.Range("O36").Value = "12.5 °C"
$ws.Range("E37").Value = "2026-02-28 17:19:37"
$ws.Range("O37").Value = "6.9 °C"
$ws.Range("E38").Value = "2026-02-28 17:19:39"
$ws.Range("K38").Value = "7.8 MJ/m2"
$ws.Range("O38").Value = "11.8 °C"
$ws.Range("E39").Value = "2026-02-28 17:19:41"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = "57%"
$fmtSource.Copy()
$ws.Range("H39").PasteSpecial(-4122)
$ws.Range("K39").Value = "14.8 MJ/m2"
$ws.Range("N39").Value = "-1.8 °C 16:59 TU"
$ws.Range("E40").Value = "2026-02-28 17:19:44"
$ws.Range("O40").Value = "7.2 °C"
$ws.Range("E41").Value = "2026-02-28 17:19:46"
$ws.Range("E42").Value = "2026-02-28 17:19:48"
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = "88%"
$fmtSource.Copy()
$ws.Range("H42").PasteSpecial(-4122)
$ws.Range("O42").Value = "11.1 °C"
$ws.Range("E43").Value = "2026-02-28 17:19:51"
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = "79%"
$fmtSource.Copy()
$ws.Range("H43").PasteSpecial(-4122)
$ws.Range("K43").Value = "6.4 MJ/m2"
$ws.Range("O43").Value = "7.1 °C"
$ws.Range("E44").Value = "2026-02-28 17:19:53"
$ws.Range("H44").NumberFormat = "@"
$ws.Range("H44").Value = "90%"
$fmtSource.Copy()
$ws.Range("H44").PasteSpecial(-4122)
$ws.Range("I44").Value = "0.1 mm"
$ws.Range("K44").Value = "13.0 MJ/m2"
$ws.Range("E45").Value = "2026-02-28 17:19:55"
$ws.Range("K45").Value = "5.7 MJ/m2"
$ws.Range("E46").Value = "2026-02-28 17:19:58"
$ws.Range("K46").Value = "4.7 MJ/m2"

$excel.CutCopyMode = 0
